$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text, matching the
# original inline-string cell type (Excel would otherwise auto-convert them to numbers).
$ws.Range('D2').Value = '37.717.88'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.026.94'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.32'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.607'
$ws.Range('E6').Value = '  -2.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.88'
$ws.Range('E7').Value = '  -2.45%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -3.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0824'
$ws.Range('E10').Value = '  +1.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.103'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '2.327.31'
$ws.Range('E12').Value = '  -2.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.43'
$ws.Range('E13').Value = '  -3.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.03'
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.766'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.18'
$ws.Range('E16').Value = '  -2.74%  '
$ws.Range('D17').Value = '2.012.35'
$ws.Range('E17').Value = '  -3.26%  '
$ws.Range('D18').Value = '37.691.27'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.44'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.89'
$ws.Range('E20').Value = '  -6.73%  '
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.49'
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.01'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.32'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('E28').Value = '  -3.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.74'
$ws.Range('E29').Value = '  -2.88%  '
$ws.Range('E30').Value = '  -5.95%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  +7.63%  '
$ws.Range('E33').Value = '  -4.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0603'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('E35').Value = '  -4.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.42'
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.29'
$ws.Range('E37').Value = '  -2.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.40'
$ws.Range('E38').Value = '  +1.24%  '
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.93'
$ws.Range('E40').Value = '  +3.16%  '
$ws.Range('D41').Value = '1.535.00'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.38'
$ws.Range('E43').Value = '  -3.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.80'
$ws.Range('E44').Value = '  -2.72%  '
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.07'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.10'
$ws.Range('E47').Value = '  -3.31%  '
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').Value = '2.217.47'
$ws.Range('E51').Value = '  -1.94%  '
